$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.867.28"
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").Value = "'1.906.91"
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("D5").Value = "'313.54"
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D7").Value = "'0.5019"
$ws.Range("E7").Value = '  +4.25%  '
$ws.Range("D8").Value = "'0.3818"
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  -1.12%  '
$ws.Range("D10").Value = "'0.9079"
$ws.Range("E10").Value = '  -2.85%  '
$ws.Range("D11").Value = "'20.83"
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = "'1.964.56"
$ws.Range("E12").Value = '  +2.38%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = "'0.07667"
$ws.Range("E13").Value = '  -1.67%  '
$ws.Range("E14").Value = '  -0.69%  '
$ws.Range("D15").Value = "'91.67"
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D17").Value = "'0.000008720"
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D19").Value = "'27.894.28"
$ws.Range("E19").Value = '  -0.91%  '
$ws.Range("D20").Value = "'14.54"
$ws.Range("E20").Value = '  -2.26%  '
$ws.Range("D21").Value = "'5.169"
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("E22").Value = '  -0.80%  '
$ws.Range("D23").Value = "'6.603"
$ws.Range("E23").Value = '  -0.81%  '
$ws.Range("D24").Value = "'154.39"
$ws.Range("D25").Value = "'1.881"
$ws.Range("E25").Value = '  -2.13%  '
$ws.Range("D26").Value = "'2.230"
$ws.Range("E26").Value = '  +5.33%  '
$ws.Range("D28").Value = "'115.30"
$ws.Range("E28").Value = '  -1.17%  '
$ws.Range("D29").Value = "'4.913"
$ws.Range("E29").Value = '  -1.27%  '
$ws.Range("D30").Value = "'0.08976"
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("E31").Value = '  -3.65%  '
$ws.Range("E32").Value = '  -1.93%  '
$ws.Range("E33").Value = '  -1.14%  '
$ws.Range("D34").Value = "'4.649"
$ws.Range("E34").Value = '  -0.92%  '
$ws.Range("D35").Value = "'0.02062"
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = "'2.552"
$ws.Range("E36").Value = '  -3.71%  '
$ws.Range("D37").Value = "'0.5591"
$ws.Range("E37").Value = '  +1.55%  '
$ws.Range("D38").Value = "'1.095"
$ws.Range("E38").Value = '  -1.37%  '
$ws.Range("E39").Value = '  +0.82%  '
$ws.Range("D40").Value = "'0.05252"
$ws.Range("E40").Value = '  -1.39%  '
$ws.Range("D41").Value = "'6.957"
$ws.Range("E41").Value = '  -1.01%  '
$ws.Range("D42").Value = "'8.499"
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  -1.02%  '
$ws.Range("D44").Value = "'111.42"
$ws.Range("E44").Value = '  +2.97%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = "'0.4818"
$ws.Range("E45").Value = '  -0.59%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'10.56"
$ws.Range("E46").Value = '  -1.94%  '
$ws.Range("D47").Value = "'1.002"
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("E48").Value = '  -1.59%  '
$ws.Range("D49").Value = "'67.53"
$ws.Range("E49").Value = '  -1.10%  '
$ws.Range("D50").Value = "'0.06068"
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("D51").Value = "'0.9014"
$ws.Range("E51").Value = '  -0.10%  '
